$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh: insert this week's two new price rows (Murcott, Primera/Segunda)
# at the top of the data block (row 69-70), pushing all older rows down by two.
$ws.Rows.Item(69).Resize(2).Insert()

# New row 69: Murcott / Primera
$ws.Cells.Item(69, 1).Value = 7
$ws.Cells.Item(69, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(69, 3).Value = "Ñuble"
$ws.Cells.Item(69, 4).Value = 44484
$ws.Cells.Item(69, 5).Value = 16
$ws.Cells.Item(69, 6).Value = "Fruta"
$ws.Cells.Item(69, 7).Value = 100102
$ws.Cells.Item(69, 8).Value = "Cítricos"
$ws.Cells.Item(69, 9).Value = 100102004
$ws.Cells.Item(69, 10).Value = "Mandarina"
$ws.Cells.Item(69, 11).Value = "Murcott"
$ws.Cells.Item(69, 12).Value = "Primera"
$ws.Cells.Item(69, 13).Value = 240
$ws.Cells.Item(69, 14).Value = 5500
$ws.Cells.Item(69, 15).Value = 6000
$ws.Cells.Item(69, 16).Value = 5750
$ws.Cells.Item(69, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(69, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(69, 19).Value = 575
$ws.Cells.Item(69, 20).Value = 10

# New row 70: Murcott / Segunda
$ws.Cells.Item(70, 1).Value = 7
$ws.Cells.Item(70, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(70, 3).Value = "Ñuble"
$ws.Cells.Item(70, 4).Value = 44484
$ws.Cells.Item(70, 5).Value = 16
$ws.Cells.Item(70, 6).Value = "Fruta"
$ws.Cells.Item(70, 7).Value = 100102
$ws.Cells.Item(70, 8).Value = "Cítricos"
$ws.Cells.Item(70, 9).Value = 100102004
$ws.Cells.Item(70, 10).Value = "Mandarina"
$ws.Cells.Item(70, 11).Value = "Murcott"
$ws.Cells.Item(70, 12).Value = "Segunda"
$ws.Cells.Item(70, 13).Value = 240
$ws.Cells.Item(70, 14).Value = 4500
$ws.Cells.Item(70, 15).Value = 5000
$ws.Cells.Item(70, 16).Value = 4750
$ws.Cells.Item(70, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(70, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(70, 19).Value = 475
$ws.Cells.Item(70, 20).Value = 10
